$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63 (shifts existing rows 63-70 down to 64-71)
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with a new price record (same as the old row 63's
# data except for the date and origin, which reflect a newer weekly entry).
$ws.Range("A63").Value = 3
$ws.Range("B63").Value = "Femacal de La Calera"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 44876
$ws.Range("E63").Value = 5
$ws.Range("F63").Value = 100112022
$ws.Range("G63").Value = "Arveja Verde"
$ws.Range("H63").Value = "Perfection"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 38
$ws.Range("K63").Value = 18000
$ws.Range("L63").Value = 18000
$ws.Range("M63").Value = 18000
$ws.Range("N63").Value = "$/malla 25 kilos"
$ws.Range("O63").Value = "Provincia de Limarí"
$ws.Range("P63").Value = 720
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
